# "vuelvo a estrategia original" — revert the "grilla de pruebas" sheet's
# scenario back to the original SELL-side strategy: flip the BUY/SELL
# direction, restore the original position-size inputs, and re-extend the
# TP ladder (rows 11-14) that had been truncated, while clearing the old
# ad-hoc SUM shortcut in E17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grilla de pruebas")

# Core scenario inputs
$ws.Range("B1").Value = "SELL"
$ws.Range("B3").Value = 2222.55
$ws.Range("F3").Value = 16.997
$ws.Range("B9").Value = 20
$ws.Range("B10").Value = 20

# Re-extend the TP ladder formulas in E/F down through row 14 (fill-down
# of the existing row-10 pattern), matching rows that had been left blank.
$ws.Range("E11").Formula = '=E10*(1+$B$6/100)'
$ws.Range("F11").Formula = '=IF($B$1="BUY",F10*(1-$B$5/100),F10*(1+$B$5/100))'
$ws.Range("E12").Formula = '=E11*(1+$B$6/100)'
$ws.Range("F12").Formula = '=IF($B$1="BUY",F11*(1-$B$5/100),F11*(1+$B$5/100))'
$ws.Range("E13").Formula = '=E12*(1+$B$6/100)'
$ws.Range("F13").Formula = '=IF($B$1="BUY",F12*(1-$B$5/100),F12*(1+$B$5/100))'
$ws.Range("E14").Formula = '=E13*(1+$B$6/100)'
$ws.Range("F14").Formula = '=IF($B$1="BUY",F13*(1-$B$5/100),F13*(1+$B$5/100))'

# Clear the old "ataque seria" shortcut formula in E17 so the cell is blank.
$ws.Range("E17").ClearContents()

# Restore the prior selection location.
$ws.Range("E14").Select()
